# Debugging, Refactoring, and more.pptx -- apply commit "Color scheme coa added"
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1 (title slide): merge the two title runs into a single run.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleTr = $slide1.Shapes.Item(1).TextFrame.TextRange
$titleAll = $titleTr.Characters(1, $titleTr.Length)
$titleAll.Text = "Debugging, Refactoring, Switch Blocks, and This"

# ---------------------------------------------------------------------------
# 2. New slide "Exercise" / "COA_color_scheme" inserted at position 9, right
#    before the existing "Debugging" slide (which -- along with everything
#    after it -- simply shifts down by one).  We build it by duplicating the
#    title slide (same ctrTitle/subTitle layout) so the placeholder/ lstStyle
#    structure matches, then move it into place and replace its text.
# ---------------------------------------------------------------------------
$dup = $slide1.Duplicate()
$newSlide = $dup.Item(1)
$newSlide.MoveTo(9)

$newSlide.Shapes.Item(1).Name = "Title 3"
$newSlide.Shapes.Item(2).Name = "Subtitle 4"

$newTitleTr = $newSlide.Shapes.Item(1).TextFrame.TextRange
$newTitleTr.Characters(1, $newTitleTr.Length).Text = "Exercise"

$newSubTr = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newSubTr.Text = "COA_color_scheme"

# ---------------------------------------------------------------------------
# 3. "CSS Refactoring" slide: merge the two runs of the "Create classes for
#    large CSS changes in JS" bullet into one run.
# ---------------------------------------------------------------------------
$cssSlide = $p.Slides.Item(7)
$cssTr = $cssSlide.Shapes.Item(2).TextFrame.TextRange
$cssPara = $cssTr.Paragraphs(7, 1)
$cssTr.Characters($cssPara.Start, $cssPara.Length).Text = "Create classes for large CSS changes in JS"

# ---------------------------------------------------------------------------
# 4. "JS Refactoring" slide: merge the two runs of the "Use functions" bullet
#    into one run.
# ---------------------------------------------------------------------------
$jsSlide = $p.Slides.Item(8)
$jsTr = $jsSlide.Shapes.Item(2).TextFrame.TextRange
$jsPara = $jsTr.Paragraphs(1, 1)
$jsTr.Characters($jsPara.Start, $jsPara.Length).Text = "Use functions"
